$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 7월 2~3주차 계획 commit: fill in the plan text for 천성호 (row 7 = C7/D7)
$ws.Range("C7").Value = "모터 기능 구현 및 webOS 개발을 위한 환경 설정(Ubuntu 설치 및 webOS emulator 설치)"
$ws.Range("D7").Value = "버튼 입력 받기, 식사 기능 구현"

# move the active selection to E7, matching the author's final cursor position
$excel.Goto($ws.Range("E7"))
